$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 0.303731
$ws.Range("G2").Value = 0.06384099999999999
$ws.Range("H2").Value = 0.02855056315031281
$ws.Range("I2").Value = 0.05595910377461311
$ws.Range("J2").Value = 0.3596901037746131
$ws.Range("K2").Value = 0.2477718962253869
